$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 block: add "n=2" label in A2 ---
$ws.Range("A2").Value = "n=2"

# --- Row 4: sign flip on F4 ---
$ws.Range("F4").Value = -1.1629537872309501

# --- Row 7 block: add "n=2" label in A7, rename header text ---
$ws.Range("A7").Value = "n=2"
$ws.Range("B7").Value = "Time-var system: p(x) = x_2"

# --- Row 9: update peak values, extend with E9/F9 ---
$ws.Range("B9").Value = 1.2500000000698901
$ws.Range("C9").Value = 1.2500000053924201
$ws.Range("D9").Value = 0.95570317361146495
$ws.Range("E9").Value = 0.913761594166
$ws.Range("F9").Value = 0.91118772095365697

# --- New block at row 11-13: Time-var system: p(x) = x_1 ---
$ws.Range("B11:F11").Merge()
$ws.Range("B11:F11").HorizontalAlignment = -4108
$ws.Range("B11").Value = "Time-var system: p(x) = x_1"
$ws.Range("A11").Value = "n=2"

$ws.Range("A12").Value = "order"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 5

$ws.Range("A13").Value = "Peak"
$ws.Range("B13").Value = 1.2500000047261901
$ws.Range("C13").Value = 1.2500000053924201
$ws.Range("D13").Value = 1.1978182268083599
$ws.Range("E13").Value = 0.854326158386573

$ws.Range("F13").Select()
